$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row number, Coin (B), Link (C), Price (D), Volume(1h) (E)
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '28.731.58', '  +6.89%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.811.43', '  +4.88%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9991', '  +0.18%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '250.94', '  +3.77%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9992', '  +0.15%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4963', '  +1.47%  '),
    @(8, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '43.00', '  +3.59%  '),
    @(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2791', '  +7.46%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06379', '  +2.61%  '),
    @(11, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.802.88', '  +4.31%  '),
    @(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '16.74', '  +4.64%  '),
    @(13, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07117', '  +3.14%  '),
    @(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6477', '  +6.41%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.712', '  +5.00%  '),
    @(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '81.87', '  +5.97%  '),
    @(17, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '28.694.24', '  +7.71%  '),
    @(18, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9989', '  +0.05%  '),
    @(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007379', '  +2.88%  '),
    @(20, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.9988', '  +0.13%  '),
    @(21, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.29', '  +7.47%  '),
    @(22, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.036.50', '  +4.21%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.627', '  +4.58%  '),
    @(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.896', '  +3.90%  '),
    @(25, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.320', '  +4.11%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '142.60', '  +2.95%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '16.04', '  +4.63%  '),
    @(28, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.881', '  +5.56%  '),
    @(29, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '112.85', '  +6.35%  '),
    @(30, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.390', '  +0.67%  '),
    @(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.174', '  +5.66%  '),
    @(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08362', '  +4.50%  '),
    @(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.832', '  +3.91%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04971', '  +9.71%  '),
    @(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.088', '  +7.95%  '),
    @(36, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6757', '  +8.11%  '),
    @(37, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.666', '  +2.81%  '),
    @(38, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.681', '  +9.12%  '),
    @(39, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9563', '  +2.26%  '),
    @(40, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.134', '  +3.61%  '),
    @(41, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01589', '  +5.77%  '),
    @(42, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.981', '  +5.87%  '),
    @(43, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9996', '  +0.26%  '),
    @(44, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '101.01', '  +1.69%  '),
    @(45, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4111', '  +6.61%  '),
    @(46, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.208', '  +4.31%  '),
    @(47, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1223', '  +5.22%  '),
    @(48, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05497', '  +1.97%  '),
    @(49, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '8.199', '  +3.15%  '),
    @(50, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '31.42', '  +4.15%  '),
    @(51, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.3618', '  +7.09%  ')
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]

    # Price (D) and Volume(1h) (E) values are stored as plain text in the workbook, but
    # several of them look like numbers (e.g. "0.9991", "43.00", "5.320") so Excel would
    # otherwise silently reinterpret them as numeric values and drop the formatting
    # (trailing zeros, thousand-dot grouping, etc). Force text entry, then clear the
    # formatting change so no stray cell style is left behind.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item[3]
    $dCell.ClearFormats()

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $item[4]
    $eCell.ClearFormats()
}
